# matrix plot rebuild with Labeled_Matrix object and hypoxanthine labeling
#
# Populates the "Originating_base" column (G) for the four ribonucleotide
# monomers, the terminal-group / tag rows and the bare-atom rows, and moves
# the view/selection to match the author's final state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Originating_base (column G) for the core A/C/G/U monomers ---
$ws.Range("G2").Value = "ADE"
$ws.Range("G3").Value = "CYT"
$ws.Range("G4").Value = "GUA"
$ws.Range("G5").Value = "URI"

# --- Originating_base for the terminal-group rows (5'OH, 5'P, 3'OH, 3'PO3, 3'cP) ---
$ws.Range("G113").Value = "none"
$ws.Range("G114").Value = "none"
$ws.Range("G115").Value = "non"
$ws.Range("G116").Value = "none"
$ws.Range("G117").Value = "none"

# --- Originating_base for the isobaric-tag rows (tag/tag2/tag3/tag4) ---
$ws.Range("G119").Value = "CYT"
$ws.Range("G120").Value = "CYT"
$ws.Range("G121").Value = "CYT"
$ws.Range("G122").Value = "CYT"

# --- Originating_base for the bare-atom rows (Hydrogen/Carbon/Oxygen/Nitrogen) ---
$ws.Range("G123").Value = "none"
$ws.Range("G124").Value = "none"
$ws.Range("G125").Value = "none"
$ws.Range("G126").Value = "none"

# --- Move the view: scroll so row 96 is at the top and select G127 ---
$ws.Range("G127").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 96
$win.ScrollColumn = 1

# --- Restore the application window geometry recorded by the author ---
$win.Left = 820
$win.Top = 760
$win.Width = 28580
$win.Height = 18360
